$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text cells to avoid Excel auto-converting numeric-looking strings
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = '@'
    $cell.Value = $text
}

Set-TextCell $ws 2 4 "29.495.46"
Set-TextCell $ws 2 5 "  +2.08%  "

Set-TextCell $ws 3 4 "1.857.37"

Set-TextCell $ws 4 4 "0.9996"
Set-TextCell $ws 4 5 "  +0.01%  "

Set-TextCell $ws 5 4 "245.54"
Set-TextCell $ws 5 5 "  +0.20%  "

Set-TextCell $ws 6 4 "0.6960"
Set-TextCell $ws 6 5 "  +0.98%  "

Set-TextCell $ws 7 5 "  +0.02%  "

Set-TextCell $ws 8 4 "0.3077"
Set-TextCell $ws 8 5 "  +0.75%  "

Set-TextCell $ws 9 4 "0.07698"
Set-TextCell $ws 9 5 "  +0.05%  "

Set-TextCell $ws 10 4 "23.64"
Set-TextCell $ws 10 5 "  +0.97%  "

Set-TextCell $ws 11 4 "0.07788"
Set-TextCell $ws 11 5 "  -0.33%  "

Set-TextCell $ws 12 4 "5.164"
Set-TextCell $ws 12 5 "  +1.41%  "

Set-TextCell $ws 13 4 "1.854.75"
Set-TextCell $ws 13 5 "  +1.20%  "

Set-TextCell $ws 14 4 "0.6944"
Set-TextCell $ws 14 5 "  +1.80%  "

Set-TextCell $ws 15 4 "91.17"
Set-TextCell $ws 15 5 "  +0.80%  "

Set-TextCell $ws 16 4 "6.341"
Set-TextCell $ws 16 5 "  -1.63%  "

Set-TextCell $ws 17 4 "29.484.49"
Set-TextCell $ws 17 5 "  +2.05%  "

Set-TextCell $ws 18 4 "0.000008320"
Set-TextCell $ws 18 5 "  +0.11%  "

Set-TextCell $ws 19 4 "2.101.14"
Set-TextCell $ws 19 5 "  +1.22%  "

Set-TextCell $ws 20 4 "238.56"
Set-TextCell $ws 20 5 "  -1.71%  "

Set-TextCell $ws 21 5 "  +0.12%  "

Set-TextCell $ws 22 4 "0.9996"
Set-TextCell $ws 22 5 "  -0.04%  "

Set-TextCell $ws 23 4 "7.619"
Set-TextCell $ws 23 5 "  +2.01%  "

Set-TextCell $ws 24 5 "  +0.04%  "

Set-TextCell $ws 25 4 "0.1496"
Set-TextCell $ws 25 5 "  +1.21%  "

Set-TextCell $ws 26 4 "160.04"
Set-TextCell $ws 26 5 "  -0.81%  "

Set-TextCell $ws 27 4 "8.899"
Set-TextCell $ws 27 5 "  +0.97%  "

Set-TextCell $ws 28 4 "18.28"
Set-TextCell $ws 28 5 "  +0.40%  "

Set-TextCell $ws 29 4 "1.533"
Set-TextCell $ws 29 5 "  -0.85%  "

Set-TextCell $ws 30 4 "4.250"
Set-TextCell $ws 30 5 "  +0.80%  "

Set-TextCell $ws 31 4 "4.149"
Set-TextCell $ws 31 5 "  -0.17%  "

Set-TextCell $ws 32 5 "  +1.89%  "

Set-TextCell $ws 33 4 "0.05108"
Set-TextCell $ws 33 5 "  -0.07%  "

Set-TextCell $ws 34 4 "0.7767"
Set-TextCell $ws 34 5 "  +1.49%  "

Set-TextCell $ws 35 4 "1.880"
Set-TextCell $ws 35 5 "  +1.99%  "

Set-TextCell $ws 36 5 "  +0.66%  "

Set-TextCell $ws 37 4 "2.689"
Set-TextCell $ws 37 5 "  -0.16%  "

Set-TextCell $ws 38 4 "1.316.00"
Set-TextCell $ws 38 5 "  +7.74%  "

Set-TextCell $ws 39 4 "0.01878"
Set-TextCell $ws 39 5 "  +1.58%  "

Set-TextCell $ws 40 5 "  +0.94%  "

Set-TextCell $ws 41 4 "0.9553"
Set-TextCell $ws 41 5 "  +0.88%  "

Set-TextCell $ws 42 4 "106.03"
Set-TextCell $ws 42 5 "  -2.33%  "

Set-TextCell $ws 43 4 "5.771"
Set-TextCell $ws 43 5 "  +0.77%  "

Set-TextCell $ws 44 5 "  +0.15%  "

Set-TextCell $ws 45 4 "9.833"
Set-TextCell $ws 45 5 "  +3.07%  "

Set-TextCell $ws 46 2 "BabyDogeCoin"
Set-TextCell $ws 46 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextCell $ws 46 4 "0.00000000125"
Set-TextCell $ws 46 5 "  +2.55%  "

Set-TextCell $ws 47 2 "RocketPoolETH"
Set-TextCell $ws 47 3 "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
Set-TextCell $ws 47 4 "2.000.70"
Set-TextCell $ws 47 5 "  +1.31%  "

Set-TextCell $ws 48 2 "Mantle"
Set-TextCell $ws 48 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell $ws 48 4 "0.5235"
Set-TextCell $ws 48 5 "  +1.39%  "

Set-TextCell $ws 49 2 "RenderToken"
Set-TextCell $ws 49 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell $ws 49 4 "1.788"
Set-TextCell $ws 49 5 "  +2.30%  "

Set-TextCell $ws 50 2 "Aave"
Set-TextCell $ws 50 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextCell $ws 50 4 "63.17"
Set-TextCell $ws 50 5 "  -1.77%  "

Set-TextCell $ws 51 2 "Aptos"
Set-TextCell $ws 51 3 "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell $ws 51 4 "6.961"
Set-TextCell $ws 51 5 "  +0.94%  "
